$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 224.33333
$ws.Range("I19").Value = 189.2
$ws.Range("K19").Value = 189.2
$ws.Range("M19").Value = -14.19999999999999
$ws.Range("H40").Value = 1451.579
$ws.Range("I40").Value = 1400
$ws.Range("J40").Value = 1498
$ws.Range("K40").Value = 1400
$ws.Range("L40").Value = 1498
$ws.Range("M40").Value = -1225
$ws.Range("N40").Value = -1848
$ws.Range("H62").Value = 2199.6155
$ws.Range("I62").Value = 2174.375
$ws.Range("J62").Value = 2240
$ws.Range("K62").Value = 2174.375
$ws.Range("L62").Value = 2240
$ws.Range("M62").Value = -1550.375
$ws.Range("N62").Value = -3488
$ws.Range("H65").Value = 2199.6155
$ws.Range("I65").Value = 2174.375
$ws.Range("J65").Value = 2240
$ws.Range("K65").Value = 10871.875
$ws.Range("L65").Value = 11200
$ws.Range("M65").Value = -7751.875
$ws.Range("N65").Value = -17440
$ws.Range("H70").Value = 51315.75
$ws.Range("I70").Value = 334433.34
$ws.Range("J70").Value = 1353.8235
$ws.Range("K70").Value = 1003300.02
$ws.Range("L70").Value = 4061.4705
$ws.Range("M70").Value = -1003030.02
$ws.Range("N70").Value = -4601.470499999999
$ws.Range("H73").Value = 51315.75
$ws.Range("I73").Value = 334433.34
$ws.Range("J73").Value = 1353.8235
$ws.Range("K73").Value = 1003300.02
$ws.Range("L73").Value = 4061.4705
$ws.Range("M73").Value = -1002364.02
$ws.Range("N73").Value = -5933.470499999999
$ws.Range("H82").Value = 10116.728
$ws.Range("I82").Value = 1128.4
$ws.Range("J82").Value = 100000
$ws.Range("K82").Value = 3385.2
$ws.Range("L82").Value = 300000
$ws.Range("M82").Value = -2979.2
$ws.Range("N82").Value = -300812
$ws.Range("H85").Value = 10116.728
$ws.Range("I85").Value = 1128.4
$ws.Range("J85").Value = 100000
$ws.Range("K85").Value = 3385.2
$ws.Range("L85").Value = 300000
$ws.Range("M85").Value = -1981.2
$ws.Range("N85").Value = -302808
$ws.Range("H98").Value = 6267.2856
$ws.Range("I98").Value = 4294.467
$ws.Range("J98").Value = 11199.333
$ws.Range("K98").Value = 4294.467
$ws.Range("L98").Value = 11199.333
$ws.Range("M98").Value = -2796.467
$ws.Range("N98").Value = -14195.333
$ws.Range("H111").Value = 1000
$ws.Range("I111").Value = 1000
$ws.Range("K111").Value = 3000
$ws.Range("M111").Value = 67
$ws.Range("H116").Value = 28574086
$ws.Range("I116").Value = 66668770
$ws.Range("J116").Value = 3075
$ws.Range("K116").Value = 66668770
$ws.Range("L116").Value = 3075
$ws.Range("M116").Value = -66665328
$ws.Range("N116").Value = -9959
$ws.Range("H122").Value = 6267.2856
$ws.Range("I122").Value = 4294.467
$ws.Range("J122").Value = 11199.333
$ws.Range("K122").Value = 12883.401
$ws.Range("L122").Value = 33597.999
$ws.Range("M122").Value = -10433.401
$ws.Range("N122").Value = -38497.999
$ws.Range("H127").Value = 958.7474999999999
$ws.Range("I127").Value = 442.85715
$ws.Range("K127").Value = 1328.57145
$ws.Range("M127").Value = 3631.42855
$ws.Range("H129").Value = 982.6
$ws.Range("J129").Value = 1096.8334
$ws.Range("L129").Value = 3290.5002
$ws.Range("N129").Value = -13290.5002
$ws.Range("H132").Value = 2791.2368
$ws.Range("I132").Value = 2092.238
$ws.Range("J132").Value = 3654.7058
$ws.Range("K132").Value = 6276.714
$ws.Range("L132").Value = 10964.1174
$ws.Range("M132").Value = -3746.714
$ws.Range("N132").Value = -16024.1174
$ws.Range("H138").Value = 2819775.8
$ws.Range("I138").Value = 6452654.5
$ws.Range("J138").Value = 4294.825
$ws.Range("K138").Value = 19357963.5
$ws.Range("L138").Value = 12884.475
$ws.Range("M138").Value = -19352823.5
$ws.Range("N138").Value = -23164.475
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1518.2759
$ws.Range("I61").Value = 1205.6957
$ws.Range("J61").Value = 2716.5
$ws.Range("K61").Value = 1205.6957
$ws.Range("L61").Value = 2716.5
$ws.Range("M61").Value = -993.6957
$ws.Range("N61").Value = -3140.5
$ws.Range("H136").Value = 1518.2759
$ws.Range("I136").Value = 1205.6957
$ws.Range("J136").Value = 2716.5
$ws.Range("K136").Value = 3617.0871
$ws.Range("L136").Value = 8149.5
$ws.Range("M136").Value = -1067.0871
$ws.Range("N136").Value = -13249.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 24000
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = $null
$ws.Range("H22").Value = 15498.25
$ws.Range("I22").Value = 20330.334
$ws.Range("J22").Value = 1002
$ws.Range("K22").Value = 20330.334
$ws.Range("L22").Value = 1002
$ws.Range("M22").Value = -20157.334
$ws.Range("N22").Value = -1348
$ws.Range("H99").Value = 1695.4584
$ws.Range("I99").Value = 1573.3334
$ws.Range("J99").Value = 1817.5834
$ws.Range("K99").Value = 1573.3334
$ws.Range("L99").Value = 1817.5834
$ws.Range("M99").Value = -75.33339999999998
$ws.Range("N99").Value = -4813.5834
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2019.881
$ws.Range("I31").Value = 1453
$ws.Range("J31").Value = 3437.0833
$ws.Range("K31").Value = 1453
$ws.Range("L31").Value = 3437.0833
$ws.Range("M31").Value = -1158
$ws.Range("N31").Value = -4027.0833
$ws.Range("H34").Value = 2019.881
$ws.Range("I34").Value = 1453
$ws.Range("J34").Value = 3437.0833
$ws.Range("K34").Value = 1453
$ws.Range("L34").Value = 3437.0833
$ws.Range("M34").Value = -1251
$ws.Range("N34").Value = -3841.0833
$ws.Range("H132").Value = 423537.1
$ws.Range("I132").Value = 644490.9
$ws.Range("J132").Value = 1716.1818
$ws.Range("K132").Value = 1933472.7
$ws.Range("L132").Value = 5148.5454
$ws.Range("M132").Value = -1930942.7
$ws.Range("N132").Value = -10208.5454
$ws.Range("H141").Value = 34251.2
$ws.Range("J141").Value = 32814
$ws.Range("L141").Value = 32814
$ws.Range("N141").Value = -43174
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1326.6
$ws.Range("I5").Value = 1637.2941
$ws.Range("J5").Value = 666.375
$ws.Range("K5").Value = 4911.8823
$ws.Range("L5").Value = 1999.125
$ws.Range("M5").Value = -4799.8823
$ws.Range("N5").Value = -2223.125
$ws.Range("H75").Value = 3448.7058
$ws.Range("I75").Value = 1228.25
$ws.Range("J75").Value = 4131.923
$ws.Range("K75").Value = 3684.75
$ws.Range("L75").Value = 12395.769
$ws.Range("M75").Value = -2686.75
$ws.Range("N75").Value = -14391.769
$ws.Range("H78").Value = 3448.7058
$ws.Range("I78").Value = 1228.25
$ws.Range("J78").Value = 4131.923
$ws.Range("K78").Value = 11054.25
$ws.Range("L78").Value = 37187.307
$ws.Range("M78").Value = -6062.25
$ws.Range("N78").Value = -47171.307
$ws.Range("H97").Value = 910
$ws.Range("I97").Value = 500
$ws.Range("J97").Value = 1183.3334
$ws.Range("K97").Value = 1500
$ws.Range("L97").Value = 3550.0002
$ws.Range("M97").Value = -1004
$ws.Range("N97").Value = -4542.0002
$ws.Range("H133").Value = 4043.5
$ws.Range("I133").Value = 1541.2858
$ws.Range("K133").Value = 4623.857400000001
$ws.Range("M133").Value = 436.1425999999992
$ws.Range("H134").Value = 4527.39
$ws.Range("I134").Value = 1760.5883
$ws.Range("J134").Value = 6487.2085
$ws.Range("K134").Value = 5281.7649
$ws.Range("L134").Value = 19461.6255
$ws.Range("M134").Value = -211.7649000000001
$ws.Range("N134").Value = -29601.6255
$ws.Range("H135").Value = 1326.6
$ws.Range("I135").Value = 1637.2941
$ws.Range("J135").Value = 666.375
$ws.Range("K135").Value = 14735.6469
$ws.Range("L135").Value = 5997.375
$ws.Range("M135").Value = -12200.6469
$ws.Range("N135").Value = -11067.375
$ws.Range("H136").Value = 4571.3335
$ws.Range("I136").Value = 1102.5
$ws.Range("K136").Value = 3307.5
$ws.Range("M136").Value = 1792.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1414.2858
$ws.Range("J46").Value = 1420
$ws.Range("L46").Value = 1420
$ws.Range("N46").Value = -1796
$ws.Range("H61").Value = 14209.3125
$ws.Range("I61").Value = 21039.1
$ws.Range("K61").Value = 21039.1
$ws.Range("M61").Value = -20837.1
$ws.Range("H113").Value = 14209.3125
$ws.Range("I113").Value = 21039.1
$ws.Range("K113").Value = 21039.1
$ws.Range("M113").Value = -18869.1
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4475
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 4950
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 4950
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -6198
$ws.Range("H65").Value = 4475
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 4950
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 24750
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -30990
$ws.Range("H81").Value = 67825.47
$ws.Range("I81").Value = 84191.836
$ws.Range("J81").Value = 2360
$ws.Range("K81").Value = 168383.672
$ws.Range("L81").Value = 4720
$ws.Range("M81").Value = -167322.672
$ws.Range("N81").Value = -6842
$ws.Range("H84").Value = 67825.47
$ws.Range("I84").Value = 84191.836
$ws.Range("J84").Value = 2360
$ws.Range("K84").Value = 841918.36
$ws.Range("L84").Value = 23600
$ws.Range("M84").Value = -836614.36
$ws.Range("N84").Value = -34208
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = $null
